$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 183
$ws.Range("C2").Value = 63.32
$ws.Range("B3").Value = 106
$ws.Range("C3").Value = 36.68
